# Update Name of Algo
# Re-run of the KNN imputation produced new values for several imputed
# cells in columns A and B of the result data sheet. Apply the updated
# values to the corresponding cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 6.907000000000001
$ws.Range("A9").Value = -21.476
$ws.Range("B9").Value = 6.456
$ws.Range("B11").Value = 5.785
$ws.Range("A13").Value = -22.08
$ws.Range("A16").Value = -21.115
$ws.Range("B16").Value = 6.272
$ws.Range("A18").Value = -21.841
$ws.Range("A20").Value = -20.623
$ws.Range("B23").Value = 7.731
$ws.Range("B24").Value = 5.645999999999999
$ws.Range("A26").Value = -21.334
$ws.Range("B26").Value = 6.364
$ws.Range("A27").Value = -21.604
$ws.Range("A29").Value = -21.516
$ws.Range("B34").Value = 7.259
$ws.Range("A35").Value = -21.638
$ws.Range("B35").Value = 5.856
$ws.Range("A36").Value = -20.929
$ws.Range("B44").Value = 5.405
$ws.Range("A45").Value = -21.175
$ws.Range("B48").Value = 5.415
$ws.Range("B49").Value = 5.905
$ws.Range("B52").Value = 5.017999999999999
$ws.Range("A55").Value = -22.124
$ws.Range("A57").Value = -22.232
$ws.Range("B66").Value = 5.202
$ws.Range("B67").Value = 5.171
$ws.Range("A69").Value = -21.291
$ws.Range("B73").Value = 6.718999999999999
$ws.Range("A76").Value = -20.392
$ws.Range("A78").Value = -20.743
$ws.Range("B78").Value = 6.742
$ws.Range("B80").Value = 8.299000000000001
$ws.Range("A82").Value = -21.81
$ws.Range("A83").Value = -21.509
$ws.Range("B91").Value = 5.509
$ws.Range("A93").Value = -21.533
$ws.Range("A97").Value = -21.684
$ws.Range("B97").Value = 5.170999999999999
$ws.Range("B99").Value = 5.3
$ws.Range("B104").Value = 7.512
